$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three 2025-12-01 earthquake records (rows 2-4); the remaining
# rows shift up so the table now spans A1:F8 instead of A1:F11.
$ws.Range("A2:F4").Delete()

# The magnitude of the last remaining record (2025-12-04 16:34:51) was
# corrected from 3.7 to 3.8. Temporarily force text formatting so the
# value is stored as text (matching the rest of the table) instead of
# becoming a numeric cell, then restore the default formatting.
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "3.8"
$ws.Range("F8").Style = "Normal"
